$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '60.821.19'
$ws.Range("D3").Value = '2.401.20'
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.36'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.81'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("E8").Value = '  +1.95%  '
$ws.Range("D9").Value = '2.407.66'
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("E10").Value = '  +0.22%  '
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.345'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.16'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.14%  '
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("D16").Value = '2.824.08'
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").Value = '60.398.77'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").Value = '2.399.55'
$ws.Range("E18").Value = '  -0.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.14'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.68'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.70'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.06'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("E25").Value = '  -2.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.65'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '574.50'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.05'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.20%  '
$ws.Range("D29").Value = '2.515.59'
$ws.Range("E29").Value = '  -1.01%  '
$ws.Range("D30").Value = '0.0₃0939'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.08'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.05%  '
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("E33").Value = '  -2.39%  '
$ws.Range("E34").Value = '  -0.21%  '
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("E36").Value = '  +3.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '152.18'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("E38").Value = '  +0.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.60'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("E42").Value = '  +8.17%  '
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.73'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.68'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("D46").Value = '0.0₆0277'
$ws.Range("E46").Value = '  -5.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.20'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.12%  '
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0509'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.32'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.51%  '
